$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 192, shifting existing rows 192:276 down to 193:277
$ws.Rows.Item(192).Insert()

# Populate the new row 192 with the latest weekly price observation
$ws.Cells.Item(192, 1).Value = 8
$ws.Cells.Item(192, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(192, 3).Value = "Coquimbo"
$ws.Cells.Item(192, 4).Value = 45134
$ws.Cells.Item(192, 5).Value = 4
$ws.Cells.Item(192, 6).Value = 100112001
$ws.Cells.Item(192, 7).Value = "Berenjena"
$ws.Cells.Item(192, 8).Value = "Sin especificar"
$ws.Cells.Item(192, 9).Value = "Primera"
$ws.Cells.Item(192, 10).Value = 400
$ws.Cells.Item(192, 11).Value = 7500
$ws.Cells.Item(192, 12).Value = 8000
$ws.Cells.Item(192, 13).Value = 7750
$ws.Cells.Item(192, 14).Value = "$/caja 50 unidades"
$ws.Cells.Item(192, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(192, 16).Value = 155
$ws.Cells.Item(192, 17).Value = 50
$ws.Cells.Item(192, 18).Value = "Hortaliza"
